$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("R$ 1664.50", "Colégio", "Gasto"),
    @("R$ 300", "Academia", "Gasto"),
    @("R$ 709.33", "Mercado", "Gasto"),
    @("R$ 1500", "Plano de saúde", "Gasto"),
    @("R$ 2000", "Aportes", "Gasto"),
    @("R$ 1855.70", "Cartão", "Gasto")
)

$row = 9
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
